{"js": "// Replace the date line and every two-digit\u00f7one-digit division prompt in\n// the table with the new values from the commit.\nconst replacements = [\n  [\"2026-02-03 Tuesday\", \"2026-02-04 Wednesday\"],\n  [\"69\u00f79=\", \"73\u00f77=\"],\n  [\"91\u00f72=\", \"66\u00f77=\"],\n  [\"61\u00f77=\", \"37\u00f76=\"],\n  [\"68\u00f79=\", \"59\u00f72=\"],\n  [\"53\u00f75=\", \"78\u00f78=\"],\n  [\"71\u00f77=\", \"29\u00f74=\"],\n  [\"32\u00f75=\", \"49\u00f72=\"],\n  [\"64\u00f77=\", \"15\u00f72=\"],\n  [\"73\u00f76=\", \"92\u00f76=\"],\n  [\"20\u00f77=\", \"32\u00f76=\"],\n  [\"27\u00f76=\", \"62\u00f78=\"],\n  [\"36\u00f76=\", \"24\u00f76=\"],\n  [\"42\u00f77=\", \"17\u00f75=\"],\n  [\"66\u00f79=\", \"20\u00f75=\"],\n  [\"51\u00f77=\", \"53\u00f74=\"],\n  [\"14\u00f74=\", \"28\u00f75=\"],\n  [\"26\u00f74=\", \"63\u00f75=\"],\n  [\"44\u00f74=\", \"12\u00f79=\"],\n  [\"76\u00f73=\", \"90\u00f72=\"],\n  [\"26\u00f79=\", \"57\u00f79=\"],\n  [\"43\u00f79=\", \"42\u00f75=\"],\n  [\"46\u00f77=\", \"86\u00f75=\"],\n  [\"62\u00f72=\", \"72\u00f79=\"],\n  [\"38\u00f74=\", \"70\u00f77=\"],\n  [\"48\u00f74=\", \"76\u00f79=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const result of results.items) {\n    result.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every two-digit\u00f7one-digit division prompt in\n# the table with the new values from the commit.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-02-03 Tuesday\", \"2026-02-04 Wednesday\"),\n    @(\"69\u00f79=\", \"73\u00f77=\"),\n    @(\"91\u00f72=\", \"66\u00f77=\"),\n    @(\"61\u00f77=\", \"37\u00f76=\"),\n    @(\"68\u00f79=\", \"59\u00f72=\"),\n    @(\"53\u00f75=\", \"78\u00f78=\"),\n    @(\"71\u00f77=\", \"29\u00f74=\"),\n    @(\"32\u00f75=\", \"49\u00f72=\"),\n    @(\"64\u00f77=\", \"15\u00f72=\"),\n    @(\"73\u00f76=\", \"92\u00f76=\"),\n    @(\"20\u00f77=\", \"32\u00f76=\"),\n    @(\"27\u00f76=\", \"62\u00f78=\"),\n    @(\"36\u00f76=\", \"24\u00f76=\"),\n    @(\"42\u00f77=\", \"17\u00f75=\"),\n    @(\"66\u00f79=\", \"20\u00f75=\"),\n    @(\"51\u00f77=\", \"53\u00f74=\"),\n    @(\"14\u00f74=\", \"28\u00f75=\"),\n    @(\"26\u00f74=\", \"63\u00f75=\"),\n    @(\"44\u00f74=\", \"12\u00f79=\"),\n    @(\"76\u00f73=\", \"90\u00f72=\"),\n    @(\"26\u00f79=\", \"57\u00f79=\"),\n    @(\"43\u00f79=\", \"42\u00f75=\"),\n    @(\"46\u00f77=\", \"86\u00f75=\"),\n    @(\"62\u00f72=\", \"72\u00f79=\"),\n    @(\"38\u00f74=\", \"70\u00f77=\"),\n    @(\"48\u00f74=\", \"76\u00f79=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
